$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Periodic refresh of the cryptocurrency price table.
#
# For the Price column (D), several new values look like plain numbers
# (e.g. "586.49"), but the source data stores every Price/Link/Coin cell
# as literal text. A bare `.Value = "586.49"` would let Excel's normal
# "looks like a number" auto-detection turn the cell into a Number, so
# for those we prefix the text with a leading apostrophe (the standard
# "force text" trick) before assigning it; values that aren't valid
# numbers (e.g. "69.201.13", "3.925.58") do not need this.

$ws.Range("D2").Value = "69.201.13"
$ws.Range("E2").Value = "  +2.17%  "
$ws.Range("D3").Value = "3.385.67"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'586.49"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").Value = "'180.02"
$ws.Range("E6").Value = "  +2.40%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("D9").Value = "'0.197"
$ws.Range("E9").Value = "  +7.45%  "
$ws.Range("D10").Value = "'0.594"
$ws.Range("E10").Value = "  +2.38%  "
$ws.Range("E11").Value = "  +3.58%  "
$ws.Range("D12").Value = "'0.0000283"
$ws.Range("E12").Value = "  +3.63%  "
$ws.Range("D13").Value = "'680.03"
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("E14").Value = "  +2.58%  "
$ws.Range("D15").Value = "3.925.58"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "69.280.41"
$ws.Range("E16").Value = "  +2.36%  "
$ws.Range("D17").Value = "3.397.44"
$ws.Range("E17").Value = "  +1.92%  "
$ws.Range("E18").Value = "  +1.61%  "
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("E20").Value = "  +2.37%  "
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D23").Value = "'17.12"
$ws.Range("E23").Value = "  +1.19%  "
$ws.Range("D24").Value = "'102.75"
$ws.Range("E24").Value = "  +1.04%  "
$ws.Range("D25").Value = "'3.92"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("E26").Value = "  +1.66%  "
$ws.Range("D27").Value = "'9.60"
$ws.Range("E27").Value = "  +1.48%  "
$ws.Range("D28").Value = "'33.91"
$ws.Range("E28").Value = "  +2.79%  "
$ws.Range("E29").Value = "  +2.33%  "
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("D32").Value = "'3.64"
$ws.Range("E32").Value = "  +10.60%  "
$ws.Range("D33").Value = "'554.69"
$ws.Range("E33").Value = "  -2.57%  "
$ws.Range("E34").Value = "  +1.02%  "
$ws.Range("D35").Value = "'58.56"
$ws.Range("E35").Value = "  +2.37%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").Value = "3.672.11"
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("D40").Value = "0.0₃0719"
$ws.Range("E40").Value = "  +7.01%  "
$ws.Range("E41").Value = "  +3.03%  "
$ws.Range("D42").Value = "'2.68"
$ws.Range("E42").Value = "  +2.47%  "
$ws.Range("E43").Value = "  +1.49%  "
$ws.Range("D46").Value = "'2.68"
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("E47").Value = "  +1.21%  "
$ws.Range("E48").Value = "  +5.60%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").Value = "'133.49"
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("D51").Value = "'2.64"
$ws.Range("E51").Value = "  +3.90%  "

# Two pairs of adjacent rows swapped rank order in this refresh:
#   row 38/39: InjectiveProtocol <-> Kaspa
#   row 44/45: VeChain <-> ApeXProtocol
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.140"
$ws.Range("E38").Value = "  +4.33%  "
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").Value = "'35.69"
$ws.Range("E39").Value = "  +1.84%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "'3.32"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0423"
$ws.Range("E45").Value = "  +3.64%  "
